$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,23
$row2[0,0] = 0.0468220186133757
$row2[0,1] = 0.00111824543683717
$row2[0,2] = 0.0340884496068105
$row2[0,3] = 0.0350263328764158
$row2[0,4] = 0.00342688117740423
$row2[0,5] = 0.00876560132746555
$row2[0,6] = 0.0455234110093067
$row2[0,7] = 0.000108217300339081
$row2[0,8] = 0.00483370608181228
$row2[0,9] = 0.0332227112040978
$row2[0,10] = 0
$row2[0,11] = 0.0518360868624197
$row2[0,12] = 0.0054469374504004
$row2[0,13] = 0.000180362167231801
$row2[0,14] = 0.95451266142414
$row2[0,15] = 0.0481927710843374
$row2[0,16] = [double]"7.21448668927206e-05"
$row2[0,17] = 0.989178269966092
$row2[0,18] = 0.000649303802034485
$row2[0,19] = 0
$row2[0,20] = 0.938460428540509
$row2[0,21] = 0.00137075247096169
$row2[0,22] = 0.00115431787028353
$ws.Range("B2:X2").Value = $row2

$row3 = New-Object 'object[,]' 1,23
$row3[0,0] = 0.889041194718996
$row3[0,1] = 0.00101002813649809
$row3[0,2] = 0.0114710338359426
$row3[0,3] = [double]"3.60724334463603e-05"
$row3[0,4] = 0.00328259144361879
$row3[0,5] = 0.0322487555010461
$row3[0,6] = 0.00721448668927206
$row3[0,7] = 0.0371546064497511
$row3[0,8] = 0.959418512372845
$row3[0,9] = 0.00854916672678739
$row3[0,10] = 0.000685376235480846
$row3[0,11] = 0
$row3[0,12] = [double]"7.21448668927206e-05"
$row3[0,13] = 0.00198398383954982
$row3[0,14] = 0.0353149123439867
$row3[0,15] = 0.00176754923887165
$row3[0,16] = 0.989935791068466
$row3[0,17] = 0.0026332876415843
$row3[0,18] = [double]"7.21448668927206e-05"
$row3[0,19] = 0.998701392395931
$row3[0,20] = 0.0125892792727797
$row3[0,21] = 0.00310222927638699
$row3[0,22] = 0.0362167231801457
$ws.Range("B3:X3").Value = $row3

$row4 = New-Object 'object[,]' 1,23
$row4[0,0] = 0.0400404011254599
$row4[0,1] = 0.0327898420027415
$row4[0,2] = 0.00306615684294062
$row4[0,3] = 0.964540797922228
$row4[0,4] = 0.992857658177621
$row4[0,5] = 0.944881321693961
$row4[0,6] = 0.00158718707163985
$row4[0,7] = 0.000144289733785441
$row4[0,8] = 0.000108217300339081
$row4[0,9] = 0.00115431787028353
$row4[0,10] = 0.000144289733785441
$row4[0,11] = 0.947262102301421
$row4[0,12] = 0.994264483082029
$row4[0,13] = 0.997799581559772
$row4[0,14] = 0.00396796767909963
$row4[0,15] = 0.945314190895318
$row4[0,16] = [double]"3.60724334463603e-05"
$row4[0,17] = 0.0028136498088161
$row4[0,18] = 0.999098189163841
$row4[0,19] = 0.00104610056994445
$row4[0,20] = 0.03885001082173
$row4[0,21] = 0.964360435754996
$row4[0,22] = 0.962520741649232
$ws.Range("B4:X4").Value = $row4

$row5 = New-Object 'object[,]' 1,23
$row5[0,0] = 0.0235913714739196
$row5[0,1] = 0.965009739557031
$row5[0,2] = 0.951157925113628
$row5[0,3] = 0.000396796767909963
$row5[0,4] = 0.000432869201356323
$row5[0,5] = 0.0135632349758315
$row5[0,6] = 0.945566697929442
$row5[0,7] = 0.962592886516124
$row5[0,8] = 0.0354952745112185
$row5[0,9] = 0.956893442031599
$row5[0,10] = 0.999170334030734
$row5[0,11] = 0.000613231368588125
$row5[0,12] = 0.000180362167231801
$row5[0,13] = 0
$row5[0,14] = 0.00591587908520309
$row5[0,15] = 0.00418440227977779
$row5[0,16] = 0.00995599163119544
$row5[0,17] = 0.00515835798282952
$row5[0,18] = 0.000180362167231801
$row5[0,19] = 0.000216434600678162
$row5[0,20] = 0.0100281364980882
$row5[0,21] = 0.0311305100642089
$row5[0,22] = [double]"7.21448668927206e-05"
$ws.Range("B5:X5").Value = $row5

Write-Output "Done"